# Corrección en parseo de Renta Fija y Venta Simultanea.
# - Para la renta fija no hay que considerar los que dicen Retrov Nominal:
#   se eliminan las filas 8-10 (BCHIAB1211, BTANN-AG, BTP0600433 - RENTA FIJA).
# - Para la venta de simultaneas no se estaba haciendo el cálculo para
#   determinar el valor de la cantidad: se actualiza la columna E (cantidad)
#   de las filas de VENTA/SIMULTANEA.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the three RENTA FIJA rows (rows 8, 9, 10) that should not be
# considered (they correspond to "Retrov Nominal" entries). Deleting
# row 8 three times removes rows 8, 9 and 10 and shifts everything below
# up by three rows.
$ws.Range("A8:J10").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# After the shift, the former rows 11, 12, 13 (VENTA / SIMULTANEA) are now
# rows 8, 9 and 10. Update their "cantidad" (column E) with the corrected
# calculation.
$ws.Range("E8").Value = 17107918.896
$ws.Range("E9").Value = 18318864
$ws.Range("E10").Value = 1213863200

$wb.Save()
